# =====================================================================
# Edit script: applies the bilibili event-list data refresh described
# by the commit (gh-pages data regeneration at 456a3b4).
#
# Sheet map (matches xl/workbook.xml order):
#   1 = 展览       (Exhibitions)
#   2 = 演出       (Performances)
#   3 = 本地生活   (Local life)
#   4 = 全部类型   (All types / combined)
# =====================================================================

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item('展览')
$ws2 = $wb.Worksheets.Item('演出')
$ws3 = $wb.Worksheets.Item('本地生活')
$ws4 = $wb.Worksheets.Item('全部类型')

function Set-Num($ws, $addr, $val) {
    $ws.Range($addr).Value = $val
}

function Set-Text($ws, $addr, $val) {
    $ws.Range($addr).Value = $val
}

# Some cells (column B, and the new row's B/E) hold strings that LOOK
# like dates (e.g. '2025-01-03'). Excel's COM Value setter auto-detects
# those and silently converts them to date serial numbers. Prefixing
# with an apostrophe forces Excel to keep them as literal text, matching
# the original workbook's inlineStr storage.
function Set-TextForceLiteral($ws, $addr, $val) {
    $ws.Range($addr).Value = "'" + $val
}

# ---------------------------------------------------------------------
# Sheet 展览 (1): 'want to go' counter (F) grew for a number of events
# ---------------------------------------------------------------------
Set-Num $ws1 "F3" 349
Set-Num $ws1 "F7" 3943
Set-Num $ws1 "F9" 796
Set-Num $ws1 "F10" 2394
Set-Num $ws1 "F12" 54
Set-Num $ws1 "F13" 239
Set-Num $ws1 "F15" 216
Set-Num $ws1 "F16" 208
Set-Num $ws1 "F17" 3676
Set-Num $ws1 "F21" 357
Set-Num $ws1 "F22" 251
Set-Num $ws1 "F23" 57

# ---------------------------------------------------------------------
# Sheet 演出 (2): same kind of counter refresh
# ---------------------------------------------------------------------
Set-Num $ws2 "F10" 104
Set-Num $ws2 "F22" 85

# ---------------------------------------------------------------------
# Sheet 本地生活 (3): counter refresh + one new event appended as row 7
# ---------------------------------------------------------------------
Set-Num $ws3 "F4" 2138
Set-Num $ws3 "F6" 27

# New row 7 (copy formatting from row 6 first, so A7 keeps the bold
# centered/bordered style used by the rest of column A)
$ws3.Range("A6").Copy()
$ws3.Range("A7").PasteSpecial(-4122)
Set-Num  $ws3 "A7" 6
Set-TextForceLiteral $ws3 "B7" '2025-01-03'
Set-Text $ws3 "C7" '广州·GAF2025插画艺术节单日门票（1/03-01/06）'
Set-Text $ws3 "D7" '新东港东路1000号 广州市保利世贸博览馆'
Set-Text $ws3 "E7" '2025.01.03 10:00-01.06 18:00'
Set-Num  $ws3 "F7" 0
Set-Num  $ws3 "G7" 68
Set-Text $ws3 "H7" 'https://show.bilibili.com/platform/detail.html?id=93909'
Set-Text $ws3 "I7" '//i1.hdslb.com/bfs/openplatform/202410/tl2B52I11729564760337.jpeg'

# ---------------------------------------------------------------------
# Sheet 全部类型 (4): counter refresh for rows outside the shifted block
# ---------------------------------------------------------------------
Set-Num $ws4 "F4" 2138
Set-Num $ws4 "F9" 349
Set-Num $ws4 "F15" 27
Set-Num $ws4 "F17" 3943
Set-Num $ws4 "F21" 104
Set-Num $ws4 "F22" 796
Set-Num $ws4 "F23" 2394
Set-Num $ws4 "F25" 54
Set-Num $ws4 "F27" 239
Set-Num $ws4 "F29" 216
Set-Num $ws4 "F30" 208
Set-Num $ws4 "F48" 85

# ---------------------------------------------------------------------
# Sheet 全部类型 (4): the '环形宇宙' event (row 33) dropped out of the feed,
# so every following event (rows 34-47) shifts up by one row, and the
# newly scraped 'GAF2025插画艺术节' event fills the vacated row 47.
# (column A, the plain row index, is untouched throughout.)
# ---------------------------------------------------------------------
# -- row 33 --
Set-TextForceLiteral $ws4 "B33" '2024-11-23'
Set-Text $ws4 "C33" '广州·cooperative kingdom综合Only同人展'
Set-Text $ws4 "D33" '西环路1号 广州岭南会展中心'
Set-Text $ws4 "E33" '2024.11.23 10:00-11.23 17:00'
Set-Num  $ws4 "F33" 325
Set-Num  $ws4 "G33" 55
Set-Text $ws4 "H33" 'https://show.bilibili.com/platform/detail.html?id=92654'
Set-Text $ws4 "I33" '//i1.hdslb.com/bfs/openplatform/202409/HnYng40s1726647875350.jpeg'
# -- row 34 --
Set-TextForceLiteral $ws4 "B34" '2024-11-24'
Set-Text $ws4 "C34" '广州·【限时早鸟8折】奇妙人声之旅 · RESOUND理想人声阿卡贝拉音乐会 '
Set-Text $ws4 "D34" '广州市二沙岛晴波路33号  星海音乐厅（交响乐演奏厅）'
Set-Text $ws4 "E34" '2024.11.24 20:00-11.24 21:30'
Set-Num  $ws4 "F34" 2
Set-Num  $ws4 "G34" 144
Set-Text $ws4 "H34" 'https://show.bilibili.com/platform/detail.html?id=90940'
Set-Text $ws4 "I34" '//i2.hdslb.com/bfs/openplatform/202408/q7p66BEy1724037045076.jpeg'
# -- row 35 --
Set-TextForceLiteral $ws4 "B35" '2024-11-27'
Set-Text $ws4 "C35" '广州·三重唱Ohashi Trio（大桥トリ才） 2024年巡演'
Set-Text $ws4 "D35" '人民北路875号（广州市少年宫内） 广州蓓蕾剧院'
Set-Text $ws4 "E35" '2024.11.27 19:30-11.27 21:00'
Set-Num  $ws4 "F35" 7
Set-Num  $ws4 "G35" 380
Set-Text $ws4 "H35" 'https://show.bilibili.com/platform/detail.html?id=91847'
Set-Text $ws4 "I35" '//i2.hdslb.com/bfs/openplatform/202409/ggAAQH8D1725369168304.jpeg'
# -- row 36 --
Set-TextForceLiteral $ws4 "B36" '2024-11-30'
Set-Text $ws4 "C36" '广州·浪潮动漫游戏展&国潮文化节'
Set-Text $ws4 "D36" '奥体南路12号 优托邦(奥体旗舰店)'
Set-Text $ws4 "E36" '2024.11.30 10:00-11.30 17:00'
Set-Num  $ws4 "F36" 238
Set-Num  $ws4 "G36" 39
Set-Text $ws4 "H36" 'https://show.bilibili.com/platform/detail.html?id=93669'
Set-Text $ws4 "I36" '//i1.hdslb.com/bfs/openplatform/202410/yVboTDM51729149825365.jpeg'
# -- row 37 --
Set-TextForceLiteral $ws4 "B37" '2024-12-06'
Set-Text $ws4 "C37" '广州·2024设计周'
Set-Text $ws4 "D37" '新港东路1000号 保利世贸博览馆'
Set-Text $ws4 "E37" '2024.12.06 09:00-12.09 17:00'
Set-Num  $ws4 "F37" 51
Set-Num  $ws4 "G37" 85
Set-Text $ws4 "H37" 'https://show.bilibili.com/platform/detail.html?id=91734'
Set-Text $ws4 "I37" '//i2.hdslb.com/bfs/openplatform/202408/PHONloTK1724306564681.jpeg'
# -- row 38 --
Set-TextForceLiteral $ws4 "B38" '2024-12-07'
Set-Text $ws4 "C38" '广州·YAYA动漫游戏嘉年华&二次元盛典'
Set-Text $ws4 "D38" '花城大道84号北门 珠江奥莱城展览中心'
Set-Text $ws4 "E38" '2024.12.07 10:00-12.07 17:00'
Set-Num  $ws4 "F38" 357
Set-Num  $ws4 "G38" 29.9
Set-Text $ws4 "H38" 'https://show.bilibili.com/platform/detail.html?id=93574'
Set-Text $ws4 "I38" '//i1.hdslb.com/bfs/openplatform/202410/jfQR0PU31728719723121.jpeg'
# -- row 39 --
Set-TextForceLiteral $ws4 "B39" '2024-12-07'
Set-Text $ws4 "C39" '广州·《型月》同人ONLY'
Set-Text $ws4 "D39" '西环路1号 广州岭南会展中心'
Set-Text $ws4 "E39" '2024.12.07 10:00-12.07 17:00'
Set-Num  $ws4 "F39" 251
Set-Num  $ws4 "G39" 60
Set-Text $ws4 "H39" 'https://show.bilibili.com/platform/detail.html?id=93092'
Set-Text $ws4 "I39" '//i0.hdslb.com/bfs/openplatform/202409/7PA42qC31727424596027.jpeg'
# -- row 40 --
Set-TextForceLiteral $ws4 "B40" '2024-12-07'
Set-Text $ws4 "C40" '广州·漫潮动漫游戏嘉年华02'
Set-Text $ws4 "D40" '东沙大道16号 广州健康方舟'
Set-Text $ws4 "E40" '2024.12.07 09:30-12.07 18:00'
Set-Num  $ws4 "F40" 57
Set-Num  $ws4 "G40" 55
Set-Text $ws4 "H40" 'https://show.bilibili.com/platform/detail.html?id=93596'
Set-Text $ws4 "I40" '//i2.hdslb.com/bfs/openplatform/202410/6TUJkmIE1729081310913.jpeg'
# -- row 41 --
Set-TextForceLiteral $ws4 "B41" '2024-12-08'
Set-Text $ws4 "C41" '广州·梁祝之父：何占豪指挥《梁祝》65周年大型东方交响音乐会'
Set-Text $ws4 "D41" '东风中路299号 广州中山纪念堂'
Set-Text $ws4 "E41" '2024.12.08 19:30-12.08 21:10'
Set-Num  $ws4 "F41" 8
Set-Num  $ws4 "G41" 70
Set-Text $ws4 "H41" 'https://show.bilibili.com/platform/detail.html?id=92833'
Set-Text $ws4 "I41" '//i1.hdslb.com/bfs/openplatform/202409/y8ck801y1726297263642.jpeg'
# -- row 42 --
Set-TextForceLiteral $ws4 "B42" '2024-12-11'
Set-Text $ws4 "C42" '广州·安田丽（安田レイ）「无形之线」2024巡演'
Set-Text $ws4 "D42" '新滘中路88号海珠同创汇东一街11号 声音共和Livehouse'
Set-Text $ws4 "E42" '2024.12.11 20:00-12.11 22:00'
Set-Num  $ws4 "F42" 46
Set-Num  $ws4 "G42" 320
Set-Text $ws4 "H42" 'https://show.bilibili.com/platform/detail.html?id=91909'
Set-Text $ws4 "I42" '//i0.hdslb.com/bfs/openplatform/202409/2821JdMa1725357077006.jpeg'
# -- row 43 --
Set-TextForceLiteral $ws4 "B43" '2024-12-14'
Set-Text $ws4 "C43" '广州·变形金刚音乐会40周年特变版'
Set-Text $ws4 "D43" '广州大道中1229号 广东艺术剧院'
Set-Text $ws4 "E43" '2024.12.14 19:30-12.14 21:30'
Set-Num  $ws4 "F43" 53
Set-Num  $ws4 "G43" 171
Set-Text $ws4 "H43" 'https://show.bilibili.com/platform/detail.html?id=90033'
Set-Text $ws4 "I43" '//i0.hdslb.com/bfs/openplatform/202407/RAV6qAVB1722168641097.jpeg'
# -- row 44 --
Set-TextForceLiteral $ws4 "B44" '2024-12-20'
Set-Text $ws4 "C44" '广州·小野丽莎2024“倾爱多彩”唱游世界音乐之旅 纪念专场'
Set-Text $ws4 "D44" '中山纪念堂 中山纪念堂'
Set-Text $ws4 "E44" '2024.12.20 20:00-12.20 22:00'
Set-Num  $ws4 "F44" 28
Set-Num  $ws4 "G44" 380
Set-Text $ws4 "H44" 'https://show.bilibili.com/platform/detail.html?id=87739'
Set-Text $ws4 "I44" '//i0.hdslb.com/bfs/openplatform/202406/HCPstM8c1718868579079.jpeg'
# -- row 45 --
Set-TextForceLiteral $ws4 "B45" '2024-12-24'
Set-Text $ws4 "C45" '广州·德国美因茨名家管弦乐团 2025 新年音乐会'
Set-Text $ws4 "D45" '人民北路875号（广州市少年宫内） 广州蓓蕾剧院'
Set-Text $ws4 "E45" '2024.12.24 19:30-12.24 21:00'
Set-Num  $ws4 "F45" 1
Set-Num  $ws4 "G45" 126
Set-Text $ws4 "H45" 'https://show.bilibili.com/platform/detail.html?id=93359'
Set-Text $ws4 "I45" '//i0.hdslb.com/bfs/openplatform/202410/HaoFdo471728632672864.jpeg'
# -- row 46 --
Set-TextForceLiteral $ws4 "B46" '2024-12-29'
Set-Text $ws4 "C46" '广州·维也纳皇家交响乐团2025新年音乐会'
Set-Text $ws4 "D46" '人民北路696号 广州友谊剧院'
Set-Text $ws4 "E46" '2024.12.29 20:00-12.30 21:45'
Set-Num  $ws4 "F46" 51
Set-Num  $ws4 "G46" 280
Set-Text $ws4 "H46" 'https://show.bilibili.com/platform/detail.html?id=89837'
Set-Text $ws4 "I46" '//i2.hdslb.com/bfs/openplatform/202407/OzlirVhz1721882951190.jpeg'
# -- row 47 --
Set-TextForceLiteral $ws4 "B47" '2025-01-03'
Set-Text $ws4 "C47" '广州·GAF2025插画艺术节单日门票（1/03-01/06）'
Set-Text $ws4 "D47" '新东港东路1000号 广州市保利世贸博览馆'
Set-Text $ws4 "E47" '2025.01.03 10:00-01.06 18:00'
Set-Num  $ws4 "F47" 0
Set-Num  $ws4 "G47" 68
Set-Text $ws4 "H47" 'https://show.bilibili.com/platform/detail.html?id=93909'
Set-Text $ws4 "I47" '//i1.hdslb.com/bfs/openplatform/202410/tl2B52I11729564760337.jpeg'

Write-Output 'edit.ps1 completed successfully'
